# Generate Report for Handoff
# Updates the localization-status report: marks files as "Ready for handoff"
# (previously "In Translation") and refreshes the handoff timestamps on the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value2 = "Ready for handoff"
$wsOverview.Range("F2").Value2 = "Ready for handoff"
$wsOverview.Range("G2").Value2 = "2016-08-20 08:47:42"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value2 = "Ready for handoff"
$wsZhCn.Range("H2").Value2 = "2016-08-20 08:47:38"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value2 = "Ready for handoff"
$wsDeDe.Range("H2").Value2 = "2016-08-20 08:47:42"

# The longer "Ready for handoff" status text widens the Status columns on
# every sheet to fit.
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
